$wb = $excel.ActiveWorkbook

# --- Sheet "body-based helpful" (sheet1) ---
$ws1 = $wb.Worksheets.Item("body-based helpful")

# Remove the "CONCERN" header from D1 (keep the cell/style, clear the value)
$ws1.Range("D1").ClearContents()

# Rename the "Which limb to move" row name to "Body Parts"
$ws1.Range("A3").Value = "Body Parts"

# Add a new "Equipment" row (moved in from the "otherwise" sheet)
$ws1.Range("A6").Value = "Equipment"
$ws1.Range("B6").Value = "Describes equipment one could use"
$ws1.Range("C6").Value = "weights, chair, box"

# --- Sheet "otherwise" (sheet3) ---
$ws3 = $wb.Worksheets.Item("otherwise")

# Remove the "Equipment" row (now duplicated on the "body-based helpful" sheet)
$ws3.Rows.Item(7).Delete()
